$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they are stored as text (matching the
# source sheet convention), not auto-converted to numbers by Excel.
$textFormatCells = @('D5', 'D6', 'D8', 'D11', 'D12', 'D13', 'D16', 'D17', 'D22', 'D23', 'D24', 'D25', 'D29', 'D30', 'D31', 'D33', 'D34', 'D35', 'D36', 'D37', 'D39', 'D41', 'D44', 'D50', 'D51')
foreach ($c in $textFormatCells) { $ws.Range($c).NumberFormat = "@" }

# Update price (D) and volume/1h change (E) columns for rows with refreshed market data
$ws.Range('D2').Value = '70.689.72'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '3.642.04'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '581.62'
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').Value = '175.62'
$ws.Range('E6').Value = '  -4.19%  '
$ws.Range('D7').Value = '3.632.46'
$ws.Range('E7').Value = '  +1.30%  '
$ws.Range('D8').Value = '0.609'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  -5.01%  '
$ws.Range('D11').Value = '6.93'
$ws.Range('E11').Value = '  +22.03%  '
$ws.Range('D12').Value = '0.608'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').Value = '48.46'
$ws.Range('E13').Value = '  -3.34%  '
$ws.Range('E14').Value = '  -2.24%  '
$ws.Range('D15').Value = '4.228.06'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').Value = '666.54'
$ws.Range('E16').Value = '  -4.39%  '
$ws.Range('D17').Value = '8.91'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '3.635.49'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D19').Value = '70.777.12'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('E21').Value = '  -2.82%  '
$ws.Range('D22').Value = '11.42'
$ws.Range('E22').Value = '  -2.58%  '
$ws.Range('D23').Value = '0.941'
$ws.Range('E23').Value = '  +1.45%  '
$ws.Range('D24').Value = '17.10'
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('D25').Value = '99.76'
$ws.Range('E25').Value = '  -4.24%  '
$ws.Range('E26').Value = '  -2.57%  '
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = '9.94'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('D30').Value = '34.67'
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('D31').Value = '3.34'
$ws.Range('E31').Value = '  -3.68%  '
$ws.Range('E32').Value = '  -0.75%  '
$ws.Range('D33').Value = '1.40'
$ws.Range('E33').Value = '  -5.58%  '
$ws.Range('D34').Value = '7.50'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = '4.02'
$ws.Range('E35').Value = '  -2.16%  '
$ws.Range('D36').Value = '585.58'
$ws.Range('E36').Value = '  -0.56%  '
$ws.Range('D37').Value = '11.08'
$ws.Range('E37').Value = '  -2.31%  '
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').Value = '58.34'
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('E40').Value = '  -0.34%  '
$ws.Range('D41').Value = '0.0457'
$ws.Range('E41').Value = '  +3.93%  '
$ws.Range('D42').Value = '3.567.15'
$ws.Range('E42').Value = '  -2.94%  '
$ws.Range('E43').Value = '  -2.67%  '
$ws.Range('D44').Value = '0.345'
$ws.Range('E44').Value = '  -0.76%  '
$ws.Range('E47').Value = '  -4.38%  '
$ws.Range('E48').Value = '  +5.40%  '
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('D50').Value = '135.91'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('D51').Value = '2.96'
$ws.Range('E51').Value = '  +1.67%  '

# Rows 45 and 46 swapped: InjectiveProtocol now ranks above PEPE, with refreshed price/volume
$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '34.58'
$ws.Range("E45").Value = '  -4.53%  '

$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D46").Value = '0.0₃0736'
$ws.Range("E46").Value = '  -5.90%  '
